$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: REPORTDATE  2020-09-30 -> 2019-09-30
$ws.Range("H2").Value = "2019-09-30 00:00:00"

# I2: BASIC_EPS  0.12 -> 0.15
$ws.Range("I2").Value = 0.15

# K2: TOTAL_OPERATE_INCOME  89942100.12 -> 70512625.65000001
$ws.Range("K2").Value = 70512625.65000001

# L2: PARENT_NETPROFIT  5357437.88 -> 6188892.97
$ws.Range("L2").Value = 6188892.97

# N2 (YSTZ), O2 (SJLTZ), P2 (BPS), Q2 (MGJYXJJE): cleared to blank
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = ""

# R2: XSMLL  30.7711064374 -> 42.051878634
$ws.Range("R2").Value = 42.051878634

# AB2: ISNEW  "1" -> "0"  (keep as text)
$ws.Range("AB2").Value = "'0"

# AC2: QDATE  "2020Q3" -> "2019Q3"
$ws.Range("AC2").Value = "2019Q3"

# AD2: DATATYPE  "2020年 三季报" -> "2019年 三季报"
$ws.Range("AD2").Value = "2019年 三季报"

# AE2: DATAYEAR  "2020" -> "2019"  (keep as text)
$ws.Range("AE2").Value = "'2019"
